$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_3_4_14"
$ws.Range("B2").Value = 0.4567298723768207
$ws.Range("C2").Value = 0.37171223060512
$ws.Range("D2").Value = 0.02392432025418956
$ws.Range("E2").Value = 0.2566516092192352
$ws.Range("F2").Value = 0.6012398600578308
$ws.Range("G2").Value = 0.3584466874599457
$ws.Range("H2").Value = 0.3441643118858337
$ws.Range("I2").Value = 0.3517256677150726
$ws.Range("A3").Value = "model_3_4_13"
$ws.Range("B3").Value = 0.4577389063096384
$ws.Range("C3").Value = 0.3769886825609916
$ws.Range("D3").Value = 0.0383354182278266
$ws.Range("E3").Value = 0.2650739331505179
$ws.Range("F3").Value = 0.6001232862472534
$ws.Range("G3").Value = 0.355436384677887
$ws.Range("H3").Value = 0.3390829563140869
$ws.Range("I3").Value = 0.3477405309677124
$ws.Range("A4").Value = "model_3_4_17"
$ws.Range("B4").Value = 0.4580330571903534
$ws.Range("C4").Value = 0.364370703337464
$ws.Range("D4").Value = 0.04781697555662257
$ws.Range("E4").Value = 0.260344400631246
$ws.Range("F4").Value = 0.5997976660728455
$ws.Range("G4").Value = 0.3626351356506348
$ws.Range("H4").Value = 0.3357397317886353
$ws.Range("I4").Value = 0.3499784171581268
$ws.Range("A5").Value = "model_3_4_15"
$ws.Range("B5").Value = 0.4581782295104782
$ws.Range("C5").Value = 0.3664524508238256
$ws.Range("D5").Value = 0.05889126693076552
$ws.Range("E5").Value = 0.2655564754080674
$ws.Range("F5").Value = 0.5996370315551758
$ws.Range("G5").Value = 0.3614474534988403
$ws.Range("H5").Value = 0.3318349421024323
$ws.Range("I5").Value = 0.3475122451782227
$ws.Range("A6").Value = "model_3_4_18"
$ws.Range("B6").Value = 0.4583182283890477
$ws.Range("C6").Value = 0.3651082871751938
$ws.Range("D6").Value = 0.05238706707219498
$ws.Range("E6").Value = 0.2624175001586527
$ws.Range("F6").Value = 0.5994821190834045
$ws.Range("G6").Value = 0.3622142970561981
$ws.Range("H6").Value = 0.3341283202171326
$ws.Range("I6").Value = 0.3489974737167358
$ws.Range("A7").Value = "model_3_4_24"
$ws.Range("B7").Value = 0.4588561908729726
$ws.Range("C7").Value = 0.3615684829484788
$ws.Range("D7").Value = 0.06215273514001007
$ws.Range("E7").Value = 0.2635827927600404
$ws.Range("F7").Value = 0.5988867878913879
$ws.Range("G7").Value = 0.3642338216304779
$ws.Range("H7").Value = 0.3306849896907806
$ws.Range("I7").Value = 0.3484461009502411
$ws.Range("A8").Value = "model_3_4_20"
$ws.Range("B8").Value = 0.4593044979292318
$ws.Range("C8").Value = 0.3665309726342261
$ws.Range("D8").Value = 0.07010863350528029
$ws.Range("E8").Value = 0.2695405060721869
$ws.Range("F8").Value = 0.5983905792236328
$ws.Range("G8").Value = 0.3614026308059692
$ws.Range("H8").Value = 0.3278797268867493
$ws.Range("I8").Value = 0.3456271290779114
$ws.Range("A9").Value = "model_3_4_19"
$ws.Range("B9").Value = 0.4595599885935927
$ws.Range("C9").Value = 0.3685925233513023
$ws.Range("D9").Value = 0.07330182954581954
$ws.Range("E9").Value = 0.2719758470462915
$ws.Range("F9").Value = 0.5981078147888184
$ws.Range("G9").Value = 0.3602265119552612
$ws.Range("H9").Value = 0.3267537951469421
$ws.Range("I9").Value = 0.3444747924804688
$ws.Range("A10").Value = "model_3_4_21"
$ws.Range("B10").Value = 0.459649421620225
$ws.Range("C10").Value = 0.3686625855539842
$ws.Range("D10").Value = 0.073586961147459
$ws.Range("E10").Value = 0.2721213034347701
$ws.Range("F10").Value = 0.5980088114738464
$ws.Range("G10").Value = 0.3601865172386169
$ws.Range("H10").Value = 0.3266532719135284
$ws.Range("I10").Value = 0.3444059789180756
$ws.Range("A11").Value = "model_3_4_23"
$ws.Range("B11").Value = 0.4596852073188684
$ws.Range("C11").Value = 0.3660962335023945
$ws.Range("D11").Value = 0.07491891274944007
$ws.Range("E11").Value = 0.2709498999097595
$ws.Range("F11").Value = 0.5979692339897156
$ws.Range("G11").Value = 0.361650675535202
$ws.Range("H11").Value = 0.3261836171150208
$ws.Range("I11").Value = 0.3449602723121643
$ws.Range("A12").Value = "model_3_4_16"
$ws.Range("B12").Value = 0.4598414981280892
$ws.Range("C12").Value = 0.3665555952801827
$ws.Range("D12").Value = 0.09108996137228176
$ws.Range("E12").Value = 0.276913914519416
$ws.Range("F12").Value = 0.5977963209152222
$ws.Range("G12").Value = 0.3613885641098022
$ws.Range("H12").Value = 0.3204816877841949
$ws.Range("I12").Value = 0.3421382904052734
$ws.Range("A13").Value = "model_3_4_22"
$ws.Range("B13").Value = 0.4599480430986721
$ws.Range("C13").Value = 0.3685486599950718
$ws.Range("D13").Value = 0.07971190740605771
$ws.Range("E13").Value = 0.2741957284588185
$ws.Range("F13").Value = 0.5976783633232117
$ws.Range("G13").Value = 0.3602515459060669
$ws.Range("H13").Value = 0.3244935870170593
$ws.Range("I13").Value = 0.3434244692325592
$ws.Range("A14").Value = "model_3_4_12"
$ws.Range("B14").Value = 0.4608451495644421
$ws.Range("C14").Value = 0.3818700365554029
$ws.Range("D14").Value = 0.1199351750973262
$ws.Range("E14").Value = 0.2968057435517027
$ws.Range("F14").Value = 0.596685528755188
$ws.Range("G14").Value = 0.3526515066623688
$ws.Range("H14").Value = 0.3103108704090118
$ws.Range("I14").Value = 0.3327262103557587
$ws.Range("A15").Value = "model_3_4_5"
$ws.Range("B15").Value = 0.4636970960529666
$ws.Range("C15").Value = 0.448860572951158
$ws.Range("D15").Value = 0.2499431937190864
$ws.Range("E15").Value = 0.3851588339145604
$ws.Range("F15").Value = 0.5935292840003967
$ws.Range("G15").Value = 0.3144325017929077
$ws.Range("H15").Value = 0.2644700407981873
$ws.Range("I15").Value = 0.2909207046031952
$ws.Range("A16").Value = "model_3_4_11"
$ws.Range("B16").Value = 0.4637999320571734
$ws.Range("C16").Value = 0.4017288223839247
$ws.Range("D16").Value = 0.1519558575273813
$ws.Range("E16").Value = 0.3207110665089121
$ws.Range("F16").Value = 0.5934154391288757
$ws.Range("G16").Value = 0.3413218259811401
$ws.Range("H16").Value = 0.299020379781723
$ws.Range("I16").Value = 0.3214150667190552
$ws.Range("A17").Value = "model_3_4_6"
$ws.Range("B17").Value = 0.4640488188825171
$ws.Range("C17").Value = 0.4382723172112011
$ws.Range("D17").Value = 0.2326280318398001
$ws.Range("E17").Value = 0.3723276325770185
$ws.Range("F17").Value = 0.5931400656700134
$ws.Range("G17").Value = 0.3204732239246368
$ws.Range("H17").Value = 0.2705753445625305
$ws.Range("I17").Value = 0.2969919741153717
$ws.Range("A18").Value = "model_3_4_9"
$ws.Range("B18").Value = 0.4640873736260201
$ws.Range("C18").Value = 0.4212781453689705
$ws.Range("D18").Value = 0.1546115925609967
$ws.Range("E18").Value = 0.3341211861841124
$ws.Range("F18").Value = 0.5930973887443542
$ws.Range("G18").Value = 0.3301686346530914
$ws.Range("H18").Value = 0.2980839908123016
$ws.Range("I18").Value = 0.3150698840618134
$ws.Range("A19").Value = "model_3_4_10"
$ws.Range("B19").Value = 0.4643265891030439
$ws.Range("C19").Value = 0.4140217921489122
$ws.Range("D19").Value = 0.1591841892987714
$ws.Range("E19").Value = 0.3310927182834928
$ws.Range("F19").Value = 0.592832624912262
$ws.Range("G19").Value = 0.3343085050582886
$ws.Range("H19").Value = 0.2964716851711273
$ws.Range("I19").Value = 0.3165028393268585
$ws.Range("A20").Value = "model_3_4_7"
$ws.Range("B20").Value = 0.4646160573074808
$ws.Range("C20").Value = 0.4391718731983731
$ws.Range("D20").Value = 0.1696262715250125
$ws.Range("E20").Value = 0.3508081816629252
$ws.Range("F20").Value = 0.5925122499465942
$ws.Range("G20").Value = 0.3199599981307983
$ws.Range("H20").Value = 0.2927897572517395
$ws.Range("I20").Value = 0.3071742057800293
$ws.Range("A21").Value = "model_3_4_8"
$ws.Range("B21").Value = 0.4649492225632076
$ws.Range("C21").Value = 0.4350990370065861
$ws.Range("D21").Value = 0.1842947909284753
$ws.Range("E21").Value = 0.3533521808414812
$ws.Range("F21").Value = 0.5921435356140137
$ws.Range("G21").Value = 0.3222836256027222
$ws.Range("H21").Value = 0.2876176536083221
$ws.Range("I21").Value = 0.3059704303741455
$ws.Range("A22").Value = "model_3_4_3"
$ws.Range("B22").Value = 0.4661589632802861
$ws.Range("C22").Value = 0.4951850277924436
$ws.Range("D22").Value = 0.3235080954002615
$ws.Range("E22").Value = 0.4405269289187165
$ws.Range("F22").Value = 0.5908046960830688
$ws.Range("G22").Value = 0.2880037724971771
$ws.Range("H22").Value = 0.2385310679674149
$ws.Range("I22").Value = 0.2647225260734558
$ws.Range("A23").Value = "model_3_4_4"
$ws.Range("B23").Value = 0.4662904317337095
$ws.Range("C23").Value = 0.4788976176970122
$ws.Range("D23").Value = 0.2780193424810777
$ws.Range("E23").Value = 0.4141781901983619
$ws.Range("F23").Value = 0.5906592011451721
$ws.Range("G23").Value = 0.2972959280014038
$ws.Range("H23").Value = 0.2545703947544098
$ws.Range("I23").Value = 0.2771897912025452
$ws.Range("A24").Value = "model_3_4_2"
$ws.Range("B24").Value = 0.4671672480960822
$ws.Range("C24").Value = 0.4987922077341442
$ws.Range("D24").Value = 0.3509205934811611
$ws.Range("E24").Value = 0.4524427545934624
$ws.Range("F24").Value = 0.5896888375282288
$ws.Range("G24").Value = 0.2859458327293396
$ws.Range("H24").Value = 0.228865385055542
$ws.Range("I24").Value = 0.2590843737125397
$ws.Range("A25").Value = "model_3_4_1"
$ws.Range("B25").Value = 0.467578884231863
$ws.Range("C25").Value = 0.5002995891793853
$ws.Range("D25").Value = 0.3598483589059108
$ws.Range("E25").Value = 0.4565355558746832
$ws.Range("F25").Value = 0.5892332792282104
$ws.Range("G25").Value = 0.2850858271121979
$ws.Range("H25").Value = 0.2257174849510193
$ws.Range("I25").Value = 0.2571478188037872
$ws.Range("A26").Value = "model_3_4_0"
$ws.Range("B26").Value = 0.4677296221541498
$ws.Range("C26").Value = 0.5170934064275281
$ws.Range("D26").Value = 0.3943330419469396
$ws.Range("E26").Value = 0.4793492558152669
$ws.Range("F26").Value = 0.5890665054321289
$ws.Range("G26").Value = 0.2755047380924225
$ws.Range("H26").Value = 0.2135581821203232
$ws.Range("I26").Value = 0.2463531792163849
